$p = $ppt.ActivePresentation

# Slide 1
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '反論'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'はんろん'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'objection, refutation, rebuttal, counterargument...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 2
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '結論'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'けつろん'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'conclusion (of an argument, discussion, study, etc.) | conclusion...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 3
$s = $p.Slides.Item(3)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '賛否'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'さんぴ'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'yes and no, for and against...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 4
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '利点'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'りてん'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'advantage, point in favor, point in favour...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 5
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'サポート'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'undefined'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'support...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 6
$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '言い換える'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'いいかえる'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'to say in other words, to put another way, to express in different words, to reword, to rephrase...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 7
$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'レベル'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'undefined'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'level, standard, amount, degree, grade, rank, class | level (plane), floor, storey (story), layer, stratum | spirit level...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 8
$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '強調'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'きょうちょう'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'emphasis, stress, highlighting, underlining, underscoring | accentuating (a feature or certain part), accenting | strong ...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 9
$s = $p.Slides.Item(9)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '部分'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'ぶぶん'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'portion, section, part...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 10
$s = $p.Slides.Item(10)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '大'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'だい'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'large, big, great, huge, vast, major, important, serious, severe | great, prominent, eminent, distinguished | -sized, as ...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 11
$s = $p.Slides.Item(11)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '付く'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'つく'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'to be attached, to be connected with, to adhere, to stick, to cling | to remain imprinted, to scar, to stain, to dye | to...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 12
$s = $p.Slides.Item(12)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '土地'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'とち'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'plot of land, lot, soil | locality, region, place...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 13
$s = $p.Slides.Item(13)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '農作物'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'のうさくぶつ'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'crops, agricultural produce...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 14
$s = $p.Slides.Item(14)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '矢張り'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'やはり'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'as expected, sure enough, just as one thought | after all (is said and done), in the end, as one would expect, in any cas...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'

# Slide 15
$s = $p.Slides.Item(15)
$s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '頼る'
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = 'たよる'
$s.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = 'to rely on, to depend on, to count on, to turn to (for help)...'
$s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = '496-510'
